$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 574, shifting existing rows 574:615 down to 575:616
$ws.Rows.Item(574).EntireRow.Insert()

# Populate the newly inserted row with the new data point.
# Force column A to text so the date-like string isn't auto-converted
# to a date serial number (it must stay a literal "yyyy/mm/dd" string,
# matching every other row in this column), then drop the now-unneeded
# number-format style so the cell matches its plain, style-less siblings.
$ws.Range("A574").NumberFormat = "@"
$ws.Range("A574").Value = "2026/01/06"
$ws.Range("A574").ClearFormats()
$ws.Range("B574").Value = "火"
$ws.Range("C574").Value = 3
$ws.Range("D574").Value = 201
